# Apply commit "feat: add 2022-Q3 data":
#  1. Insert a new worksheet named "2022-Q3" right after "总计" (i.e. before the
#     current "2022-Q2" sheet), populated with the new quarter's fund holdings.
#  2. Prepend a matching summary row to the "总计" sheet and renumber the
#     running index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q3" sheet.
#    Duplicate the existing "2022-Q2" sheet (Worksheets.Item(2)) so the new
#    sheet inherits identical column widths / header styling / cell formats,
#    then overwrite its data with the 2022-Q3 figures. Copy() with a single
#    "before" argument places the clone immediately in front of the source
#    sheet, i.e. right after "总计" - exactly where the diff wants it.
$template = $wb.Worksheets.Item(2)
$template.Copy($template)
$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Fund rows: code, name, size, total stock position, position ratio, held value, rank
$funds = @(
    @("013554", "信澳远见价值混合A", "0.92", "48.39", "2.03", "0.0187", 9),
    @("013555", "信澳远见价值混合C", "0.64", "48.39", "2.03", "0.0130", 9),
    @("090011", "大成核心双动力混合", "0.24", "92.56", "2.49", "0.0060", 9)
)

for ($i = 0; $i -lt $funds.Length; $i++) {
    $row = $i + 2
    $q3.Cells.Item($row, 1).Value = $i
    # Numeric-looking text must keep a leading apostrophe so it is stored as
    # text (matching the source data) instead of being coerced to a number.
    $q3.Cells.Item($row, 2).Value = "'" + $funds[$i][0]
    $q3.Cells.Item($row, 3).Value = $funds[$i][1]
    $q3.Cells.Item($row, 4).Value = "'" + $funds[$i][2]
    $q3.Cells.Item($row, 5).Value = "'" + $funds[$i][3]
    $q3.Cells.Item($row, 6).Value = "'" + $funds[$i][4]
    $q3.Cells.Item($row, 7).Value = "'" + $funds[$i][5]
    $q3.Cells.Item($row, 8).Value = $funds[$i][6]
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) "总计" sheet: add the 2022-Q3 summary row on top, shift the rest down.
$zj = $wb.Worksheets.Item(1)

# The sheet is about to gain a 9th data row (row 9) that didn't exist before;
# clone the formatting of the last existing row's index cell (A8) onto it so
# the running-index column keeps its styling.
$zj.Cells.Item(8, 1).Copy()
$zj.Cells.Item(9, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totals = @(
    @("2022-Q3", 3, 0.04),
    @("2022-Q2", 3, 0.02),
    @("2022-Q1", 2, 0),
    @("2021-Q4", 2, 0.4),
    @("2021-Q3", 4, 1.78),
    @("2021-Q2", 6, 0.03),
    @("2021-Q1", 4, 2.86),
    @("2020-Q4", 26, 2.86)
)

for ($i = 0; $i -lt $totals.Length; $i++) {
    $row = $i + 2
    $zj.Cells.Item($row, 1).Value = $i
    $zj.Cells.Item($row, 2).Value = $totals[$i][0]
    $zj.Cells.Item($row, 3).Value = $totals[$i][1]
    $zj.Cells.Item($row, 4).Value = $totals[$i][2]
}
